# Availability_List.xlsx update
# - Fix Row 2 / Row 3 availability windows
# - Set Availability: refresh row 5 (BUSY slot) and append new application
#   outcome rows (6-13) with BUSY/AVAILABLE slots for D001
# - Drop the stale custom column widths picked up from a previous session

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Correct the existing availability rows (2-4)
# ---------------------------------------------------------------------

# Row 2: D001 30/10 09:00-17:00 -> 09:00-10:00
$ws.Range("D2").Value = 0.41666666666666669

# Row 3: D001 31/10 10:00-15:00 -> 30/10 10:00-23:00
$ws.Range("B3").Value = 45595
$ws.Range("D3").Value = 0.95833333333333337

# Row 4 (D002) is unchanged in value/format.

# ---------------------------------------------------------------------
# Row 5: existing BUSY row gets bumped to the new application outcome
# (date + times updated, status kept as BUSY)
# ---------------------------------------------------------------------
$ws.Range("B5").Value = 45599
$ws.Range("C5").Value = "12:00 pm"
$ws.Range("D5").Value = "01:00 pm"
$ws.Range("E5").Value = "BUSY"

# ---------------------------------------------------------------------
# Insert 8 new application-outcome rows below row 5 (new rows 6-13)
# ---------------------------------------------------------------------
$ws.Rows(6).Resize(8).Insert()

$newRows = @(
    @{ Row=6;  Date=45599; Start="01:00 pm"; End="02:00 pm"; Status="BUSY" },
    @{ Row=7;  Date=45600; Start="02:00 pm"; End="03:00 pm"; Status="AVAILABLE" },
    @{ Row=8;  Date=45600; Start="03:00 pm"; End="04:00 pm"; Status="AVAILABLE" },
    @{ Row=9;  Date=45601; Start="05:00 pm"; End="06:00 pm"; Status="AVAILABLE" },
    @{ Row=10; Date=45601; Start="06:00 pm"; End="07:00 pm"; Status="AVAILABLE" },
    @{ Row=11; Date=45599; Start="03:00 pm"; End="04:00 pm"; Status="AVAILABLE" },
    @{ Row=12; Date=45599; Start="04:00 pm"; End="05:00 pm"; Status="AVAILABLE" },
    @{ Row=13; Date=45599; Start="02:00 pm"; End="03:00 pm"; Status="AVAILABLE" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = "D001"
    $ws.Range("B$row").Value = $r.Date
    $ws.Range("C$row").Value = $r.Start
    $ws.Range("D$row").Value = $r.End
    $ws.Range("E$row").Value = $r.Status
}

# ---------------------------------------------------------------------
# Drop the stale custom widths on columns B:D (matches the diff dropping
# the <cols> overrides) and re-apply the existing number formats that
# ClearFormats() would otherwise wipe off the data rows.
# ---------------------------------------------------------------------
$ws.Columns("B:D").ClearFormats()

$ws.Range("B2:B4").NumberFormat = "mm-dd-yy"
$ws.Range("C2:D4").NumberFormat = "h:mm AM/PM"
$ws.Range("B5:B13").NumberFormat = "dd/MM/yy"

# ClearFormats() on the whole column also strips the bordered/bold header
# look from B1:D1 - restore it (A1/E1 were untouched and keep theirs).
$ws.Range("B1:D1").Font.Bold = $true
$ws.Range("B1:D1").BorderAround(1, 2, 64)

# ---------------------------------------------------------------------
# Move the active selection to C6, matching the saved sheet view
# ---------------------------------------------------------------------
$ws.Range("C6").Select()
